# "upgrade left table until javakheti"
# Shuakhevi.xlsx: rename the sheet, mark a few more cells in the
# Urban/Rural breakdown as confidential ("…"), and drop the blank
# spacer row above the footnote so it moves from row 9 to row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The tab was generically named "1" - give it the municipality name.
$ws.Name = "Shuakhevi"

# Row 6 (Urban) becomes fully confidential/unavailable across the
# columns that previously had real counts.
$ws.Range("D6").Value = "…"
$ws.Range("E6").Value = "…"
$ws.Range("G6").Value = "…"
$ws.Range("I6").Value = "…"

# Row 7 (Rural): one more column (2015) flips to confidential.
$ws.Range("G7").Value = "…"

# Remove the empty row sitting between the data table and the
# footnote so the note moves up from row 9 to row 8.
$ws.Rows.Item(8).Delete()
